$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 170
$ws.Range("I8").Value = 170
$ws.Range("K8").Value = 510
$ws.Range("M8").Value = -371

$ws.Range("H41").Value = 1219.6316
$ws.Range("I41").Value = 961.9167
$ws.Range("J41").Value = 1661.4286
$ws.Range("K41").Value = 961.9167
$ws.Range("L41").Value = 1661.4286
$ws.Range("M41").Value = -521.9167
$ws.Range("N41").Value = -2541.4286

$ws.Range("H76").Value = 3245.9583
$ws.Range("I76").Value = 3186.1428
$ws.Range("J76").Value = 3329.7
$ws.Range("K76").Value = 3186.1428
$ws.Range("L76").Value = 3329.7
$ws.Range("M76").Value = -2871.1428
$ws.Range("N76").Value = -3959.7

$ws.Range("H79").Value = 3245.9583
$ws.Range("I79").Value = 3186.1428
$ws.Range("J79").Value = 3329.7
$ws.Range("K79").Value = 3186.1428
$ws.Range("L79").Value = 3329.7
$ws.Range("M79").Value = -2094.1428
$ws.Range("N79").Value = -5513.7

$ws.Range("H86").Value = 2593.5217
$ws.Range("I86").Value = 2234.1428
$ws.Range("J86").Value = 3152.5557
$ws.Range("K86").Value = 2234.1428
$ws.Range("L86").Value = 3152.5557
$ws.Range("M86").Value = -1111.1428
$ws.Range("N86").Value = -5398.5557

$ws.Range("H89").Value = 2593.5217
$ws.Range("I89").Value = 2234.1428
$ws.Range("J89").Value = 3152.5557
$ws.Range("K89").Value = 11170.714
$ws.Range("L89").Value = 15762.7785
$ws.Range("M89").Value = -5554.714
$ws.Range("N89").Value = -26994.7785

$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459

$ws.Range("H106").Value = 2800.7856
$ws.Range("I106").Value = 2600.4167
$ws.Range("J106").Value = 4003
$ws.Range("K106").Value = 2600.4167
$ws.Range("L106").Value = 4003
$ws.Range("M106").Value = -1969.4167
$ws.Range("N106").Value = -5265

$ws.Range("H127").Value = 1115.8889
$ws.Range("I127").Value = 779
$ws.Range("J127").Value = 1284.3334
$ws.Range("K127").Value = 2337
$ws.Range("L127").Value = 3853.0002
$ws.Range("M127").Value = 2623
$ws.Range("N127").Value = -13773.0002

$ws.Range("H129").Value = 853.64514
$ws.Range("I129").Value = 437.8
$ws.Range("J129").Value = 933.61536
$ws.Range("K129").Value = 1313.4
$ws.Range("L129").Value = 2800.84608
$ws.Range("M129").Value = 3686.6
$ws.Range("N129").Value = -12800.84608

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 5999
$ws.Range("I11").Value = 5999
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 5999
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -5855

$ws.Range("H102").Value = 83334830
$ws.Range("I102").Value = 83334830
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 83334830
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -83333208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2473
$ws.Range("I20").Value = 2319.6
$ws.Range("J20").Value = 2779.8
$ws.Range("K20").Value = 2319.6
$ws.Range("L20").Value = 2779.8
$ws.Range("M20").Value = -2072.6
$ws.Range("N20").Value = -3273.8

$ws.Range("H105").Value = 142859380
$ws.Range("I105").Value = 250002000
$ws.Range("K105").Value = 250002000
$ws.Range("M105").Value = -250000253

$ws.Range("H107").Value = 2059.4614
$ws.Range("I107").Value = 1385
$ws.Range("J107").Value = 2637.5715
$ws.Range("K107").Value = 1385
$ws.Range("L107").Value = 2637.5715
$ws.Range("M107").Value = 535
$ws.Range("N107").Value = -6477.5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2250
$ws.Range("I2").Value = 425
$ws.Range("J2").Value = 3466.6667
$ws.Range("K2").Value = 425
$ws.Range("L2").Value = 3466.6667
$ws.Range("M2").Value = -312
$ws.Range("N2").Value = -3692.6667

$ws.Range("H31").Value = 2232.4517
$ws.Range("I31").Value = 1050.7273
$ws.Range("J31").Value = 2882.4
$ws.Range("K31").Value = 1050.7273
$ws.Range("L31").Value = 2882.4
$ws.Range("M31").Value = -755.7273
$ws.Range("N31").Value = -3472.4

$ws.Range("H34").Value = 2232.4517
$ws.Range("I34").Value = 1050.7273
$ws.Range("J34").Value = 2882.4
$ws.Range("K34").Value = 1050.7273
$ws.Range("L34").Value = 2882.4
$ws.Range("M34").Value = -848.7273
$ws.Range("N34").Value = -3286.4

$ws.Range("H93").Value = 17528.572
$ws.Range("I93").Value = 3540
$ws.Range("J93").Value = 52500
$ws.Range("K93").Value = 3540
$ws.Range("L93").Value = 52500
$ws.Range("M93").Value = -1668
$ws.Range("N93").Value = -56244

$ws.Range("H99").Value = 2633172
$ws.Range("I99").Value = 3290977.5
$ws.Range("J99").Value = 1950
$ws.Range("K99").Value = 3290977.5
$ws.Range("L99").Value = 1950
$ws.Range("M99").Value = -3289479.5
$ws.Range("N99").Value = -4946

$ws.Range("H105").Value = 1000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("M105").Value = 1000
$ws.Range("N105").Value = -4494

$ws.Range("H126").Value = 2633172
$ws.Range("I126").Value = 3290977.5
$ws.Range("J126").Value = 1950
$ws.Range("K126").Value = 9872932.5
$ws.Range("L126").Value = 5850
$ws.Range("M126").Value = -9870462.5
$ws.Range("N126").Value = -10790

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 177.88889
$ws.Range("I2").Value = 59.75
$ws.Range("J2").Value = 272.4
$ws.Range("K2").Value = 358.5
$ws.Range("L2").Value = 1634.4
$ws.Range("M2").Value = -245.5
$ws.Range("N2").Value = -1860.4

$ws.Range("H23").Value = 398
$ws.Range("I23").Value = 1100
$ws.Range("J23").Value = 222.5
$ws.Range("K23").Value = 3300
$ws.Range("L23").Value = 667.5
$ws.Range("M23").Value = -3065
$ws.Range("N23").Value = -1137.5

$ws.Range("H33").Value = 293.18182
$ws.Range("J33").Value = 302.66666
$ws.Range("L33").Value = 1815.99996
$ws.Range("N33").Value = -2381.99996

$ws.Range("H87").Value = 2019
$ws.Range("I87").Value = 604.6667
$ws.Range("K87").Value = 1814.0001
$ws.Range("M87").Value = -566.0001

$ws.Range("H90").Value = 2019
$ws.Range("I90").Value = 604.6667
$ws.Range("K90").Value = 5442.0003
$ws.Range("M90").Value = 797.9997000000003

$ws.Range("H126").Value = 5109.6772
$ws.Range("I126").Value = 2480
$ws.Range("J126").Value = 5615.385
$ws.Range("K126").Value = 7440
$ws.Range("L126").Value = 16846.155
$ws.Range("M126").Value = -2500
$ws.Range("N126").Value = -26726.155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 32149466
$ws.Range("I70").Value = 62505252
$ws.Range("J70").Value = 20007150
$ws.Range("K70").Value = 62505252
$ws.Range("L70").Value = 20007150
$ws.Range("M70").Value = -62504982
$ws.Range("N70").Value = -20007690

$ws.Range("H73").Value = 32149466
$ws.Range("I73").Value = 62505252
$ws.Range("J73").Value = 20007150
$ws.Range("K73").Value = 62505252
$ws.Range("L73").Value = 20007150
$ws.Range("M73").Value = -62504316
$ws.Range("N73").Value = -20009022

$ws.Range("H97").Value = 626.625
$ws.Range("I97").Value = 656
$ws.Range("K97").Value = 656
$ws.Range("M97").Value = -160

$ws.Range("H126").Value = 2263.35
$ws.Range("I126").Value = 1822
$ws.Range("K126").Value = 5466
$ws.Range("M126").Value = -2996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 917.5
$ws.Range("I16").Value = 894.25
$ws.Range("J16").Value = 1150
$ws.Range("K16").Value = 894.25
$ws.Range("L16").Value = 1150
$ws.Range("M16").Value = -724.25
$ws.Range("N16").Value = -1490

$ws.Range("H22").Value = 1248.7059
$ws.Range("I22").Value = 1254.4546
$ws.Range("J22").Value = 1238.1666
$ws.Range("K22").Value = 1254.4546
$ws.Range("L22").Value = 1238.1666
$ws.Range("M22").Value = -959.4546
$ws.Range("N22").Value = -1828.1666

$ws.Range("H27").Value = 1248.7059
$ws.Range("I27").Value = 1254.4546
$ws.Range("J27").Value = 1238.1666
$ws.Range("K27").Value = 1254.4546
$ws.Range("L27").Value = 1238.1666
$ws.Range("M27").Value = -1147.4546
$ws.Range("N27").Value = -1452.1666

$ws.Range("H40").Value = 3265.5833
$ws.Range("I40").Value = 2954.5715
$ws.Range("J40").Value = 3701
$ws.Range("K40").Value = 2954.5715
$ws.Range("L40").Value = 3701
$ws.Range("M40").Value = -2818.5715
$ws.Range("N40").Value = -3973

$ws.Range("H122").Value = 14914529
$ws.Range("I122").Value = 23613096
$ws.Range("K122").Value = 70839288
$ws.Range("M122").Value = -70836838

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18524274
$ws.Range("I62").Value = 19236054
$ws.Range("K62").Value = 19236054
$ws.Range("M62").Value = -19235430

$ws.Range("H65").Value = 18524274
$ws.Range("I65").Value = 19236054
$ws.Range("K65").Value = 96180270
$ws.Range("M65").Value = -96177150

$ws.Range("H81").Value = 4314
$ws.Range("I81").Value = 678.2
$ws.Range("J81").Value = 5383.353
$ws.Range("K81").Value = 1356.4
$ws.Range("L81").Value = 10766.706
$ws.Range("M81").Value = -295.4000000000001
$ws.Range("N81").Value = -12888.706

$ws.Range("H84").Value = 4314
$ws.Range("I84").Value = 678.2
$ws.Range("J84").Value = 5383.353
$ws.Range("K84").Value = 6782
$ws.Range("L84").Value = 53833.53
$ws.Range("M84").Value = -1478
$ws.Range("N84").Value = -64441.53

$ws.Range("H113").Value = 399.3846
$ws.Range("I113").Value = 310.41177
$ws.Range("J113").Value = 567.44446
$ws.Range("K113").Value = 931.23531
$ws.Range("L113").Value = 1702.33338
$ws.Range("M113").Value = 1238.76469
$ws.Range("N113").Value = -6042.33338

$ws.Range("H126").Value = 222223000
$ws.Range("I126").Value = 222223000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 666669000
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -666666530
